$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output $ws.Range("F11").Value
Write-Output $ws.Range("B8").Value
Write-Output $ws.Range("J10").Hyperlinks.Count
